$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.322284
$ws.Range("H2").Value = 6.966852
$ws.Range("I2").Value = 0.05374183487824914
$ws.Range("J2").Value = 0.05374183487824913
$ws.Range("M2").Value = 2.027115333333333
$ws.Range("N2").Value = 6.081346
$ws.Range("O2").Value = 0.006596284565418616
$ws.Range("P2").Value = 0.006596284565418615
$ws.Range("Q2").Value = 4.707537504754667
$ws.Range("R2").Value = 42.367837542792
$ws.Range("S2").Value = 0.0003544964359246706
$ws.Range("T2").Value = 0.0003544964359246705
$ws.Range("G3").Value = 2.322284
$ws.Range("H3").Value = 6.966852
$ws.Range("I3").Value = 0.05374183487824914
$ws.Range("J3").Value = 0.05374183487824913
$ws.Range("O3").Value = 0.8344762556643375
$ws.Range("P3").Value = 0.8344762556643374
$ws.Range("Q3").Value = 595.5365071667147
$ws.Range("R3").Value = 5359.828564500432
$ws.Range("S3").Value = 0.04484628514173243
$ws.Range("T3").Value = 0.04484628514173242
$ws.Range("G4").Value = 2.322284
$ws.Range("H4").Value = 6.966852
$ws.Range("I4").Value = 0.05374183487824914
$ws.Range("J4").Value = 0.05374183487824913
$ws.Range("M4").Value = 48.84026566666667
$ws.Range("N4").Value = 146.520797
$ws.Range("O4").Value = 0.158927459770244
$ws.Range("P4").Value = 0.158927459770244
$ws.Range("Q4").Value = 113.4209675134494
$ws.Range("R4").Value = 1020.788707621044
$ws.Range("S4").Value = 0.008541053300592037
$ws.Range("T4").Value = 0.008541053300592032
$ws.Range("I5").Value = 0.5740922721009293
$ws.Range("J5").Value = 0.5740922721009293
$ws.Range("M5").Value = 2.027115333333333
$ws.Range("N5").Value = 6.081346
$ws.Range("O5").Value = 0.006596284565418616
$ws.Range("P5").Value = 0.006596284565418615
$ws.Range("Q5").Value = 50.287842017816
$ws.Range("R5").Value = 452.590578160344
$ws.Range("S5").Value = 0.003786875993585464
$ws.Range("T5").Value = 0.003786875993585463
$ws.Range("I6").Value = 0.5740922721009293
$ws.Range("J6").Value = 0.5740922721009293
$ws.Range("O6").Value = 0.8344762556643375
$ws.Range("P6").Value = 0.8344762556643374
$ws.Range("S6").Value = 0.4790663696286154
$ws.Range("T6").Value = 0.4790663696286154
$ws.Range("I7").Value = 0.5740922721009293
$ws.Range("J7").Value = 0.5740922721009293
$ws.Range("M7").Value = 48.84026566666667
$ws.Range("N7").Value = 146.520797
$ws.Range("O7").Value = 0.158927459770244
$ws.Range("P7").Value = 0.158927459770244
$ws.Range("Q7").Value = 1211.609188469212
$ws.Range("R7").Value = 10904.48269622291
$ws.Range("S7").Value = 0.09123902647872843
$ws.Range("T7").Value = 0.09123902647872839
$ws.Range("G8").Value = 16.081976
$ws.Range("H8").Value = 48.245928
$ws.Range("I8").Value = 0.3721658930208215
$ws.Range("J8").Value = 0.3721658930208215
$ws.Range("M8").Value = 2.027115333333333
$ws.Range("N8").Value = 6.081346
$ws.Range("O8").Value = 0.006596284565418616
$ws.Range("P8").Value = 0.006596284565418615
$ws.Range("Q8").Value = 32.60002013989867
$ws.Range("R8").Value = 293.400181259088
$ws.Range("S8").Value = 0.00245491213590848
$ws.Range("T8").Value = 0.00245491213590848
$ws.Range("G9").Value = 16.081976
$ws.Range("H9").Value = 48.245928
$ws.Range("I9").Value = 0.3721658930208215
$ws.Range("J9").Value = 0.3721658930208215
$ws.Range("O9").Value = 0.8344762556643375
$ws.Range("P9").Value = 0.8344762556643374
$ws.Range("Q9").Value = 4124.131163707339
$ws.Range("R9").Value = 37117.18047336605
$ws.Range("S9").Value = 0.3105636008939895
$ws.Range("T9").Value = 0.3105636008939894
$ws.Range("G10").Value = 16.081976
$ws.Range("H10").Value = 48.245928
$ws.Range("I10").Value = 0.3721658930208215
$ws.Range("J10").Value = 0.3721658930208215
$ws.Range("M10").Value = 48.84026566666667
$ws.Range("N10").Value = 146.520797
$ws.Range("O10").Value = 0.158927459770244
$ws.Range("P10").Value = 0.158927459770244
$ws.Range("Q10").Value = 785.4479802849575
$ws.Range("R10").Value = 7069.031822564617
$ws.Range("S10").Value = 0.05914737999092355
$ws.Range("T10").Value = 0.05914737999092353
